# Update the cryptos worksheet with the latest scraped values.
# Columns: A=rank index, B=Coin, C=Link, D=Price, E=Volume(1h)
# D (Price) cells are stored as literal text (e.g. "67.818.37", "0.0000288"),
# so force text format before assigning the value to avoid Excel silently
# converting the numeric-looking strings into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Rows whose Coin/Link stayed the same, only Price (D) and/or Volume (E) changed ---

Set-PriceText "D2" "68.118.17"
$ws.Range("E2").Value = "  -5.46%  "

Set-PriceText "D3" "3.715.21"
$ws.Range("E3").Value = "  -4.47%  "

$ws.Range("E4").Value = "  -0.07%  "

Set-PriceText "D5" "580.34"
$ws.Range("E5").Value = "  -2.16%  "

Set-PriceText "D6" "182.10"
$ws.Range("E6").Value = "  +8.21%  "

Set-PriceText "D7" "3.711.72"
$ws.Range("E7").Value = "  -4.34%  "

Set-PriceText "D8" "0.627"
$ws.Range("E8").Value = "  -6.96%  "

Set-PriceText "D9" "0.997"
$ws.Range("E9").Value = "  -0.31%  "

Set-PriceText "D10" "0.712"
$ws.Range("E10").Value = "  -7.40%  "

Set-PriceText "D11" "0.163"
$ws.Range("E11").Value = "  -10.20%  "

Set-PriceText "D12" "53.28"
$ws.Range("E12").Value = "  -2.63%  "

Set-PriceText "D13" "0.0000291"
$ws.Range("E13").Value = "  -10.54%  "

Set-PriceText "D14" "10.49"
$ws.Range("E14").Value = "  -7.41%  "

Set-PriceText "D15" "4.267.86"
$ws.Range("E15").Value = "  -5.34%  "

Set-PriceText "D16" "3.696.95"
$ws.Range("E16").Value = "  -5.68%  "

Set-PriceText "D17" "19.42"
$ws.Range("E17").Value = "  -7.85%  "

$ws.Range("E18").Value = "  -2.77%  "

Set-PriceText "D19" "12.87"
$ws.Range("E19").Value = "  -7.87%  "

Set-PriceText "D20" "1.12"
$ws.Range("E20").Value = "  -7.43%  "

Set-PriceText "D21" "67.870.45"
$ws.Range("E21").Value = "  -5.66%  "

Set-PriceText "D22" "407.56"
$ws.Range("E22").Value = "  -7.55%  "

Set-PriceText "D23" "4.46"
$ws.Range("E23").Value = "  -5.14%  "

Set-PriceText "D24" "88.23"
$ws.Range("E24").Value = "  -6.53%  "

Set-PriceText "D25" "3.06"
$ws.Range("E25").Value = "  -6.95%  "

Set-PriceText "D26" "12.81"
$ws.Range("E26").Value = "  -7.81%  "

Set-PriceText "D27" "10.92"
$ws.Range("E27").Value = "  -1.81%  "

# --- Rows 28/29: LEO and Toncoin swap places (with updated Price/Volume) ---

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-PriceText "D28" "3.80"
$ws.Range("E28").Value = "  -9.93%  "

$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-PriceText "D29" "6.05"
$ws.Range("E29").Value = "  +2.15%  "

Set-PriceText "D30" "9.52"
$ws.Range("E30").Value = "  -6.26%  "

Set-PriceText "D31" "32.59"
$ws.Range("E31").Value = "  -7.43%  "

Set-PriceText "D32" "7.52"
$ws.Range("E32").Value = "  -3.07%  "

Set-PriceText "D33" "12.50"
$ws.Range("E33").Value = "  -8.33%  "

# --- Rows 34-37: Bittensor, Hedera, OKB, InjectiveProtocol re-ranked ---

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-PriceText "D34" "0.117"
$ws.Range("E34").Value = "  -7.62%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-PriceText "D35" "43.57"
$ws.Range("E35").Value = "  -13.41%  "

$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-PriceText "D36" "603.42"
$ws.Range("E36").Value = "  -1.97%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-PriceText "D37" "64.73"
$ws.Range("E37").Value = "  -4.14%  "

Set-PriceText "D38" "0.0₃0894"
$ws.Range("E38").Value = "  -8.44%  "

$ws.Range("E39").Value = "  +0.35%  "

Set-PriceText "D40" "0.399"
$ws.Range("E40").Value = "  -5.20%  "

Set-PriceText "D41" "0.996"
$ws.Range("E41").Value = "  -0.45%  "

Set-PriceText "D42" "0.136"
$ws.Range("E42").Value = "  -5.63%  "

Set-PriceText "D43" "2.78"
$ws.Range("E43").Value = "  +7.11%  "

Set-PriceText "D44" "3.01"
$ws.Range("E44").Value = "  -9.66%  "

# --- Rows 45/46: VeChain and dogwifhat swap places (with updated Price/Volume) ---

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-PriceText "D45" "2.97"
$ws.Range("E45").Value = "  -7.07%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-PriceText "D46" "0.0436"
$ws.Range("E46").Value = "  -7.41%  "

Set-PriceText "D47" "9.24"
$ws.Range("E47").Value = "  -10.28%  "

Set-PriceText "D48" "2.815.44"
$ws.Range("E48").Value = "  -0.24%  "

Set-PriceText "D49" "0.134"
$ws.Range("E49").Value = "  -7.38%  "

Set-PriceText "D50" "2.70"
$ws.Range("E50").Value = "  -5.37%  "

Set-PriceText "D51" "3.06"
$ws.Range("E51").Value = "  -7.71%  "
